# Update odds values in Sheet1 to reflect the latest FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 6
$ws.Range("X2").Value = 9.5
$ws.Range("AC2").Value = 6
$ws.Range("AO2").Value = 15
$ws.Range("AW2").Value = 5

# Row 6
$ws.Range("Q6").Value = 1.53
$ws.Range("R6").Value = 2.4

# Row 7
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 10
$ws.Range("Q7").Value = 2.08
$ws.Range("R7").Value = 1.73

# Row 8
$ws.Range("G8").Value = 2.3
$ws.Range("I8").Value = 3.1
$ws.Range("J8").Value = 3.1
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 8
$ws.Range("Q8").Value = 2.3
$ws.Range("R8").Value = 1.6
$ws.Range("X8").Value = 10
$ws.Range("Y8").Value = 10
$ws.Range("AG8").Value = 8
$ws.Range("AN8").Value = 4.33
$ws.Range("AR8").Value = 81

# Row 10
$ws.Range("G10").Value = 2.55
$ws.Range("I10").Value = 3
$ws.Range("K10").Value = 2
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 7.5
$ws.Range("U10").Value = 1.95
$ws.Range("V10").Value = 1.8
$ws.Range("AC10").Value = 7.5
$ws.Range("AI10").Value = 11
$ws.Range("AJ10").Value = 29
$ws.Range("BA10").Value = 81

# Row 16
$ws.Range("M16").Value = 1.11
$ws.Range("N16").Value = 6.5

# Row 19
$ws.Range("G19").Value = 2.2
$ws.Range("I19").Value = 3.9
$ws.Range("J19").Value = 3.1
$ws.Range("AX19").Value = 26

# Row 20
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = 8
$ws.Range("J20").Value = 17.5
$ws.Range("L20").Value = 1.28
$ws.Range("O20").Value = 1.06
$ws.Range("P20").Value = 7.1
$ws.Range("R20").Value = 3.75
$ws.Range("T20").Value = 4.9
$ws.Range("X20").Value = 600
$ws.Range("Y20").Value = 120
$ws.Range("AB20").Value = 300
$ws.Range("AC20").Value = 29
$ws.Range("AD20").Value = 25
$ws.Range("AE20").Value = 45
$ws.Range("AF20").Value = 175
$ws.Range("AG20").Value = 13.5
$ws.Range("AI20").Value = 14
$ws.Range("AN20").Value = 27
$ws.Range("AO20").Value = 200
$ws.Range("AP20").Value = 90
$ws.Range("AT20").Value = 4.9
$ws.Range("AU20").Value = 11.5
